$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header row text to include Arabic translations
$ws.Range("A1").Value = "Trip Reference*   رقم الرحلة"
$ws.Range("B1").Value = "Point Reference*   رقم نقطة التوصيل"
$ws.Range("C1").Value = "Point Type* نوع نقطة التوصيل"
$ws.Range("D1").Value = "Facility*  المنطقة"
$ws.Range("E1").Value = "Agent*  الوكيل"

# Remove sample data rows (rows 2-4), keeping only the header
$ws.Range("A2:E4").Delete()

# Resize columns
$ws.Columns.Item(1).ColumnWidth = 42.5703125
$ws.Columns.Item(2).ColumnWidth = 37.5703125
$ws.Columns.Item(3).ColumnWidth = 40.5703125
$ws.Columns.Item(4).ColumnWidth = 29.5703125
$ws.Columns.Item(5).ColumnWidth = 37

# Update selection to match target state
$ws.Range("A4:XFD4").Select()
